# Updates column B ("word") for rows 2-193 with the new working set of sequences.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'wellen',
    'arten',
    'wehren',
    'gelten',
    'drohen',
    'platzen',
    'wehtun',
    'boxen',
    'liefern',
    'helfen',
    'feuern',
    'machen',
    'jubeln',
    'streichen',
    'spüren',
    'klettern',
    'stehlen',
    'leeren',
    'leiden',
    'tropfen',
    'passen',
    'saufen',
    'dienen',
    'geben',
    'regeln',
    'fällen',
    'ändern',
    'küssen',
    'mögen',
    'irren',
    'quälen',
    'lockern',
    'narren',
    'greifen',
    'klingen',
    'schämen',
    'wachsen',
    'scheiden',
    'gründen',
    'fangen',
    'wahren',
    'flehen',
    'reiten',
    'zielen',
    'räumen',
    'fallen',
    'läuten',
    'laufen',
    'schlucken',
    'sorgen',
    'lösen',
    'kennen',
    'schreiten',
    'kümmern',
    'rücken',
    'bauen',
    'kehren',
    'albern',
    'reizen',
    'wüten',
    'ruhen',
    'pflanzen',
    'münzen',
    'heben',
    'kranken',
    'stechen',
    'scheitern',
    'bitten',
    'stoppen',
    'schwören',
    'jagen',
    'hauen',
    'beißen',
    'tauschen',
    'trauen',
    'wenden',
    'fischen',
    'leisten',
    'rufen',
    'proben',
    'fragen',
    'sichern',
    'siegen',
    'lehnen',
    'schulden',
    'leugnen',
    'tollen',
    'liegen',
    'graben',
    'werfen',
    'schenken',
    'trennen',
    'meistern',
    'hacken',
    'lesen',
    'dauern',
    'zahlen',
    'töten',
    'backen',
    'bilden',
    'rasen',
    'heilen',
    'stören',
    'spielen',
    'runden',
    'baden',
    'bleiben',
    'duschen',
    'seufzen',
    'husten',
    'warnen',
    'gleichen',
    'schmecken',
    'fordern',
    'kosten',
    'scheinen',
    'nennen',
    'zeigen',
    'malen',
    'wehen',
    'lächeln',
    'altern',
    'werden',
    'betteln',
    'süßen',
    'reisen',
    'äußern',
    'weichen',
    'spinnen',
    'kichern',
    'schwächen',
    'zeichnen',
    'knurren',
    'drehen',
    'fahren',
    'pfeifen',
    'stürzen',
    'wirken',
    'stammen',
    'heulen',
    'schütteln',
    'enden',
    'strahlen',
    'hören',
    'grüßen',
    'wundern',
    'decken',
    'sperren',
    'sterben',
    'spenden',
    'bergen',
    'bellen',
    'folgen',
    'planen',
    'drucken',
    'sinken',
    'treiben',
    'streifen',
    'stecken',
    'feiern',
    'ärgern',
    'freuen',
    'fließen',
    'zählen',
    'filmen',
    'flüchten',
    'füllen',
    'schwingen',
    'tragen',
    'ehren',
    'bluten',
    'brauchen',
    'sparen',
    'loben',
    'formen',
    'lügen',
    'mauern',
    'wagen',
    'biegen',
    'sprengen',
    'achten',
    'erben',
    'gnaden',
    'zünden',
    'schrecken',
    'stärken',
    'führen',
    'trösten',
    'suchen',
    'knarren',
    'dringen',
    'zögern'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}

